$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new URL column (keeps the same style as the other headers)
$ws.Range("E1").Value = "URL"

# URLs for each tournament row, in sheet order (rows 2-13)
$urls = @(
    "https://www.espn.com/golf/leaderboard?tournamentId=401056527",
    "https://www.espn.com/golf/leaderboard/_/tournamentId/401056550",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056551",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056552",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056554",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056556",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056558",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056548",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056547",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056544",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056543",
    "https://www.espn.com/golf/leaderboard?tournamentId=401056542"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $urls[$i]
    $cell.Style = "Normal"
}

$ws.Range("D19").Select()
